$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Note: several "Price" (column D) values look like plain numbers
# (e.g. "554.90", "0.100"). Assigning those via .Value would make Excel
# auto-convert them to numeric cells and drop the exact text formatting
# (trailing zeros, etc.), which would not match the source data (kept as
# text). For those cells we instead set .Formula with a leading apostrophe
# (the same thing Excel does when a user types '554.90 into a cell), which
# forces a text cell, and then reset .Style back to 'Normal' so no stray
# "quote prefix" cell style lingers on the cell.

$ws.Range('D2').Value = '58.823.82'
$ws.Range('E2').Value = '  -1.40%  '
$ws.Range('D3').Value = '2.595.91'
$ws.Range('E3').Value = '  -2.71%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Formula = '''554.90'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.21%  '
$ws.Range('D6').Formula = '''143.02'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.38%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Formula = '''0.599'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.35%  '
$ws.Range('D9').Formula = '''6.77'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.37%  '
$ws.Range('D10').Formula = '''0.100'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E11').Value = '  +5.22%  '
$ws.Range('E12').Value = '  -1.19%  '
$ws.Range('D13').Value = '3.054.94'
$ws.Range('E13').Value = '  -2.21%  '
$ws.Range('D14').Value = '58.803.38'
$ws.Range('E14').Value = '  -1.29%  '
$ws.Range('D15').Formula = '''20.81'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.85%  '
$ws.Range('D16').Value = '2.606.63'
$ws.Range('E16').Value = '  -2.57%  '
$ws.Range('E17').Value = '  -2.33%  '
$ws.Range('E18').Value = '  -0.11%  '
$ws.Range('D19').Formula = '''337.43'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.79%  '
$ws.Range('D20').Formula = '''10.09'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -2.65%  '
$ws.Range('D21').Formula = '''6.16'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.77%  '
$ws.Range('E22').Value = '  -0.19%  '
$ws.Range('D23').Formula = '''66.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').Formula = '''0.429'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.50%  '
$ws.Range('D25').Formula = '''0.995'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.40%  '
$ws.Range('D26').Formula = '''0.161'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.36%  '
$ws.Range('D27').Formula = '''7.15'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.92%  '
$ws.Range('E28').Value = '  +0.44%  '
$ws.Range('E29').Value = '  -0.11%  '
$ws.Range('E30').Value = '  +1.49%  '
$ws.Range('D31').Formula = '''6.02'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.26%  '
$ws.Range('D32').Formula = '''154.47'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.31%  '
$ws.Range('D33').Formula = '''18.99'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.37%  '
$ws.Range('E34').Value = '  -2.16%  '
$ws.Range('D35').Formula = '''0.896'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.76%  '
$ws.Range('E36').Value = '  -0.90%  '
$ws.Range('D37').Formula = '''36.88'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.11%  '
$ws.Range('D38').Formula = '''0.857'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.09%  '
$ws.Range('D39').Formula = '''1.48'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.99%  '
$ws.Range('E40').Value = '  -0.46%  '
$ws.Range('D41').Formula = '''282.97'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.45%  '
$ws.Range('D42').Formula = '''0.997'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.25%  '
$ws.Range('D43').Formula = '''0.600'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -1.35%  '
$ws.Range('D44').Formula = '''0.0954'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.82%  '
$ws.Range('B45').Value = 'Hedera'
$ws.Range('C45').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D45').Formula = '''0.0534'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.25%  '
$ws.Range('B46').Value = 'WhiteBITCoin'
$ws.Range('C46').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D46').Formula = '''10.61'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.25%  '
$ws.Range('E47').Value = '  -0.18%  '
$ws.Range('D48').Value = '1.939.69'
$ws.Range('E48').Value = '  -1.48%  '
$ws.Range('D49').Formula = '''117.93'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.81%  '
$ws.Range('D50').Formula = '''17.96'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.58%  '
$ws.Range('E51').Value = '  -2.93%  '
